$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.163737893104553
$ws.Range("B1").Value = 2.660921335220337
$ws.Range("C1").Value = 4.235062122344971
$ws.Range("D1").Value = 3.415743350982666
$ws.Range("E1").Value = 1.2118079662323
